$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.27
$ws.Range("B8").Value = 5.792
$ws.Range("B10").Value = 6.944
$ws.Range("B12").Value = 6.444
$ws.Range("C13").Value = -12.729
$ws.Range("B18").Value = 6.873
$ws.Range("D20").Value = -8.222
$ws.Range("B25").Value = 6.629
